# ADD results from server
# Update investment-cost result values on sheets "2025", "2030" and "2035"
# to reflect the latest server run.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 2062.1809032
$ws.Range("E2").Value = 291069.2862091724
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 148652.5872276
$ws.Range("L2").Value = 509125.9821312751
$ws.Range("M2").Value = 112470.9127927
$ws.Range("N2").Value = 72560.46740536761
$ws.Range("O2").Value = 68708.80120585454

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 48847.97250432047
$ws.Range("E2").Value = 275288.7038538232
$ws.Range("I2").Value = 284097.1166367626
$ws.Range("L2").Value = 257107.5985254353
$ws.Range("M2").Value = 105992.8640084325
$ws.Range("N2").Value = 33892.28451258693
$ws.Range("O2").Value = 35046.15750099967

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 22680.65032447391
$ws.Range("B2").Value = 19944.27970611305
$ws.Range("E2").Value = 105109.2055170748
$ws.Range("I2").Value = 90040.7107790998
$ws.Range("M2").Value = 57074.83554118505
$ws.Range("N2").Value = 51192.95312943371
$ws.Range("O2").Value = 49027.54443095748
